$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# "running all Profile xlsx" - flip the Runmode column (D) from "N" to "Y"
# for every test case row so the whole suite executes.
$ws.Range("D3:D41").Value = "Y"

# Matches the resulting Excel selection after performing the above edit.
$ws.Range("D2:D41").Select()
